$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0,
    33.23912340812262,
    -8.116932572908174,
    -18.86411409169933,
    -14.2165379802521,
    -11.40996652814376,
    -16.52778345848755,
    -16.52778345848755,
    -16.52778345848766,
    -16.52778345848755
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $val = $values[$i]
    $ws.Range("E$row").Value = $val
    $ws.Range("F$row").Value = $val
}
